$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("1:1").Insert()

$headers = @("Determined","Collaborator","Family","Optimistic","Conqueror","Peacful zionist","Nurturing","Compassionate","Educator","Army Commander","Visionary","N/A","Achiever","Fiction Character","Problem Sover")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 2 + $i
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}
$ws.Range("B1:P1").Font.Color = 0

$ws.Range("B5").Value = "significant"
$ws.Range("B7").Value = "topic"
$ws.Range("C7").Value = "effect"

$data = @(
    @(2, 19.56),
    @(3, 149.79),
    @(4, 49.66),
    @(5, 137.63),
    @(6, 85.39),
    @(7, 9.6),
    @(9, 44.15),
    @(10, 85.19)
)
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 8 + $i
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
}

Write-Host "done"
